$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Base")

# The two Salesforce record Ids that were rotated for this regression sheet.
$oldCompanyId = "001q000000irAHS"
$newCompanyId = "0010x000004TMut"
$oldUserId    = "005q0000003GGfP"
$newUserId    = "0051I000000guzk"

for ($col = 1; $col -le 10; $col++) {
    $cell = $ws.Cells.Item(2, $col)
    $text = $cell.Value2
    if ($text -ne $null) {
        $newText = $text -replace [regex]::Escape($oldCompanyId), $newCompanyId
        $newText = $newText -replace [regex]::Escape($oldUserId), $newUserId
        if ($newText -ne $text) {
            if ($col -eq 10) {
                # J2 ("Username") stores the bare Id and Excel keeps it quote-prefixed
                # (style carries quotePrefix="1") because it looks numeric-ish.
                # Re-enter it with a leading apostrophe so that formatting is retained.
                $cell.Value = "'" + $newText
            } else {
                $cell.Value2 = $newText
            }
        }
    }
}

# Reflect the reviewer's final on-screen selection/scroll position for the sheet.
$ws.Range("H16").Select()
$excel.ActiveWindow.ScrollColumn = 6
